$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Re-categorize a few existing entries (ring column, B) from Adopt to Trial
$ws.Range("B2").Value = "Trial"
$ws.Range("B13").Value = "Trial"
$ws.Range("B15").Value = "Trial"

# Add two new Dependency Injection entries at the bottom of the table
$ws.Range("A76").Value = "Dagger"
$ws.Range("B76").Value = "Hold"
$ws.Range("C76").Value = "languages-and-frameworks"
$ws.Range("D76").Value = "FALSE"

$ws.Range("A77").Value = "Koin"
$ws.Range("B77").Value = "Assess"
$ws.Range("C77").Value = "languages-and-frameworks"
$ws.Range("D77").Value = "FALSE"
